$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 18 and 19 (route data shortened / restructured)
$ws.Rows("18:19").Delete() | Out-Null

# Adjust column widths: E shrinks (8 -> 7), new column F added (18)
$ws.Columns("E").ColumnWidth = 6.1666666666667
$ws.Columns("F").ColumnWidth = 17.1666666666667

# Copy formatting from column E into new column F (header + data rows) so styles line up (s=2 header / s=3 data)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("E2:E17").Copy() | Out-Null
$ws.Range("F2:F17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 1
$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "Current Node"
$ws.Range("C1").Value = "Next Node"
$ws.Range("D1").Value = "Used Capacity (leaving the node)"
$ws.Range("E1").Value = "Route"
$ws.Range("F1").Value = "Node Type"

# Row 2
$ws.Range("A2").Value = 393
$ws.Range("B2").Value = "MC"
$ws.Range("C2").Formula = "=""1"""
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 11
$ws.Range("F2").Value = "start of route"

# Row 3
$ws.Range("A3").Value = 398
$ws.Range("B3").Formula = "=""1"""
$ws.Range("C3").Value = "MCd"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 11
$ws.Range("F3").Value = "type 1 drop off"

# Row 4
$ws.Range("A4").Value = 403
$ws.Range("B4").Value = "MCd"
$ws.Range("C4").Value = "n.a."
$ws.Range("D4").Value = "n.a."
$ws.Range("E4").Value = 11
$ws.Range("F4").Value = "end of route"

# Row 5
$ws.Range("A5").Value = 403
$ws.Range("B5").Value = "MC"
$ws.Range("C5").Formula = "=""4"""
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 14
$ws.Range("F5").Value = "start of route"

# Row 6
$ws.Range("A6").Value = 407
$ws.Range("B6").Formula = "=""4"""
$ws.Range("C6").Value = "1p"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 14
$ws.Range("F6").Value = "patient pickup"

# Row 7
$ws.Range("A7").Value = 428
$ws.Range("B7").Value = "1p"
$ws.Range("C7").Value = "MCd"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 14
$ws.Range("F7").Value = "type 1 pick up"

# Row 8
$ws.Range("A8").Value = 700
$ws.Range("B8").Value = "MCd"
$ws.Range("C8").Value = "n.a."
$ws.Range("D8").Value = "n.a."
$ws.Range("E8").Value = 14
$ws.Range("F8").Value = "end of route"

# Row 9
$ws.Range("A9").Value = 861
$ws.Range("B9").Value = "MC"
$ws.Range("C9").Formula = "=""3"""
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = "start of route"

# Row 10
$ws.Range("A10").Value = 865
$ws.Range("B10").Formula = "=""3"""
$ws.Range("C10").Value = "3p"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = "type 2 drop off"

# Row 11
$ws.Range("A11").Value = 895
$ws.Range("B11").Value = "3p"
$ws.Range("C11").Formula = "=""2"""
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = "type 2 pick up"

# Row 12
$ws.Range("A12").Value = 900
$ws.Range("B12").Formula = "=""2"""
$ws.Range("C12").Value = "MCd"
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = "type 2 drop off"

# Row 13
$ws.Range("A13").Value = 903
$ws.Range("B13").Value = "MCd"
$ws.Range("C13").Value = "n.a."
$ws.Range("D13").Value = "n.a."
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = "end of route"

# Row 14
$ws.Range("A14").Value = 903
$ws.Range("B14").Value = "MC"
$ws.Range("C14").Value = "4p"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 22
$ws.Range("F14").Value = "start of route"

# Row 15
$ws.Range("A15").Value = 907
$ws.Range("B15").Value = "4p"
$ws.Range("C15").Value = "2p"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 22
$ws.Range("F15").Value = "patient drop off"

# Row 16
$ws.Range("A16").Value = 990
$ws.Range("B16").Value = "2p"
$ws.Range("C16").Value = "MCd"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 22
$ws.Range("F16").Value = "type 2 pick up"

# Row 17
$ws.Range("A17").Value = 993
$ws.Range("B17").Value = "MCd"
$ws.Range("C17").Value = "n.a."
$ws.Range("D17").Value = "n.a."
$ws.Range("E17").Value = 22
$ws.Range("F17").Value = "end of route"

# Convert the text formulas above into plain static text values, keeping original cell style
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4163) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4163) | Out-Null
$ws.Range("C5").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4163) | Out-Null
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4163) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C9").PasteSpecial(-4163) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4163) | Out-Null
$ws.Range("C11").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4163) | Out-Null
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
